# Update cryptos list with latest price/volume data
# (Updated cryptos list on Thu Sep 28 09:42:36 UTC 2023 with GitHub Actions)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a literal text value to a cell without letting Excel
# auto-coerce numeric-looking strings (e.g. "63.79", "5.90", "0.0501")
# into floating point numbers. Temporarily mark the cell as Text so the
# value round-trips byte-for-byte, then restore the cell's original
# (default) style so no stray formatting is left behind.
function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$ws.Range("D2").Value = '26.393.61'
$ws.Range("E2").Value = '  +0.39%  '
$ws.Range("D3").Value = '1.618.17'
$ws.Range("E3").Value = '  +1.24%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("E5").Value = '  -0.12%  '
$ws.Range("E6").Value = '  -0.31%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("E8").Value = '  +0.20%  '
$ws.Range("E9").Value = '  +0.09%  '
$ws.Range("E10").Value = '  +0.54%  '
$ws.Range("E11").Value = '  -0.61%  '
$ws.Range("D12").Value = '1.845.84'
$ws.Range("E12").Value = '  +1.42%  '
$ws.Range("D13").Value = '1.613.50'
$ws.Range("E13").Value = '  +0.73%  '
$ws.Range("E14").Value = '  +0.09%  '
$ws.Range("E15").Value = '  -0.20%  '
Set-TextValue "D16" '63.79'
$ws.Range("E16").Value = '  -0.33%  '
Set-TextValue "D17" '235.99'
$ws.Range("E17").Value = '  +9.12%  '
$ws.Range("D18").Value = '26.408.73'
$ws.Range("E18").Value = '  +0.49%  '
Set-TextValue "D19" '7.80'
$ws.Range("E19").Value = '  +4.34%  '
$ws.Range("E20").Value = '  +0.24%  '
$ws.Range("E21").Value = '  +0.01%  '
$ws.Range("E22").Value = '  -0.56%  '
$ws.Range("E23").Value = '  +1.03%  '
$ws.Range("E24").Value = '  +2.75%  '
Set-TextValue "D25" '147.08'
$ws.Range("E25").Value = '  +1.58%  '
$ws.Range("E27").Value = '  +0.83%  '
$ws.Range("E28").Value = '  +0.29%  '
Set-TextValue "D29" '15.54'
$ws.Range("E30").Value = '  +0.14%  '
Set-TextValue "D31" '1.15'
$ws.Range("E31").Value = '  -0.58%  '
$ws.Range("D32").Value = '1.516.67'
$ws.Range("E32").Value = '  +5.82%  '
$ws.Range("E33").Value = '  +1.29%  '
Set-TextValue "D34" '2.97'
$ws.Range("E34").Value = '  -0.09%  '
Set-TextValue "D35" '1.51'
$ws.Range("E35").Value = '  +3.09%  '
$ws.Range("E36").Value = '  +0.45%  '
Set-TextValue "D37" '0.567'
$ws.Range("E37").Value = '  +1.51%  '
$ws.Range("E38").Value = '  +0.52%  '
Set-TextValue "D39" '0.834'
$ws.Range("E39").Value = '  +0.67%  '
Set-TextValue "D40" '5.90'
$ws.Range("E40").Value = '  +2.22%  '
$ws.Range("E41").Value = '  -0.03%  '
$ws.Range("E42").Value = '  +1.42%  '
$ws.Range("D43").Value = '1.757.48'
$ws.Range("E43").Value = '  +1.37%  '
Set-TextValue "D44" '0.761'
$ws.Range("E44").Value = '  +0.11%  '
Set-TextValue "D45" '61.85'
$ws.Range("E45").Value = '  +1.39%  '
Set-TextValue "D46" '0.907'
$ws.Range("E46").Value = '  +0.28%  '
Set-TextValue "D47" '90.14'
$ws.Range("E47").Value = '  +3.45%  '
Set-TextValue "D48" '1.51'
$ws.Range("E48").Value = '  +1.70%  '
Set-TextValue "D49" '0.0501'
$ws.Range("E49").Value = '  +0.03%  '
$ws.Range("E50").Value = '  +0.96%  '
Set-TextValue "D51" '7.49'
$ws.Range("E51").Value = '  +0.75%  '
